$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Psilo")

# Insert 4 new columns before column F (old F shifts to J, G->K, H->L, I->M)
$ws.Range("F1:I1").EntireColumn.Insert()

# New column headers
$ws.Range("F1").Value = "FearNeut"
$ws.Range("G1").Value = "SadNeut"
$ws.Range("H1").Value = "HappyNeut"
$ws.Range("I1").Value = "AngryNeut"

# New data for columns F:I (FearNeut, SadNeut, HappyNeut, AngryNeut), rows 2-40
$fiData = @(
    @(0.0338983050847458, 0.0, 0.0, 0.0),
    @(0.0169491525423729, 0.05, 0.0, 0.0),
    @(0.62711864406779705, 0.125, 0.0, 0.025),
    @(0.101694915254237, 0.075, 0.0, 0.10000000000000001),
    @(0.0169491525423729, 0.10000000000000001, 0.0, 0.125),
    @(0.76271186440677996, 0.05, 0.125, 0.14999999999999999),
    @(0.0508474576271187, 0.05, 0.025, 0.025),
    @(0.0847457627118644, 0.05, 0.05, 0.05),
    @(0.28813559322033899, 0.025, 0.025, 0.0),
    @(0.0508474576271187, 0.05, 0.10000000000000001, 0.05),
    @(0.0, 0.05, 0.0, 0.025),
    @(0.0677966101694915, 0.05, 0.025, 0.0),
    @(0.0847457627118644, 0.05, 0.025, 0.22500000000000001),
    @(0.22033898305084701, 0.025, 0.125, 0.10000000000000001),
    @(0.93220338983050799, 0.075, 0.05, 0.14999999999999999),
    @(0.69491525423728795, 0.5, 0.52500000000000002, 0.52500000000000002),
    @(0.11864406779661001, 0.025, 0.0, 0.025),
    @(0.54237288135593198, 0.05, 0.025, 0.025),
    @(0.0847457627118644, 0.17499999999999999, 0.0, 0.05),
    @(0.50847457627118597, 0.025, 0.025, 0.05),
    @(0.50847457627118597, 0.025, 0.0, 0.025),
    @(0.66101694915254205, 0.05, 0.17499999999999999, 0.125),
    @(0.71186440677966101, 0.025, 0.05, 0.0),
    @(0.66101694915254205, 0.05, 0.0, 0.025),
    @(0.66101694915254205, 0.10000000000000001, 0.67500000000000004, 0.075),
    @(0.0677966101694915, 0.05, 0.025, 0.025),
    @(0.0508474576271187, 0.075, 0.05, 0.10000000000000001),
    @(0.72881355932203395, 0.10000000000000001, 0.125, 0.075),
    @(0.77966101694915302, 0.10000000000000001, 0.0, 0.10000000000000001),
    @(0.677966101694915, 0.0, 0.0, 0.0),
    @(0.0338983050847458, 0.0, 0.0, 0.0),
    @(0.0338983050847458, 0.025, 0.0, 0.025),
    @(0.50847457627118597, 0.05, 0.025, 0.05),
    @(0.0169491525423729, 0.0, 0.025, 0.05),
    @(0.677966101694915, 0.0, 0.0, 0.0),
    @(0.94915254237288105, 0.025, 0.075, 0.10000000000000001),
    @(0.677966101694915, 0.025, 0.0, 0.0),
    @(0.169491525423729, 0.14999999999999999, 0.10000000000000001, 0.10000000000000001),
    @(0.0847457627118644, 0.025, 0.025, 0.05)
)

for ($i = 0; $i -lt $fiData.Count; $i++) {
    $r = $i + 2
    $vals = $fiData[$i]
    $ws.Cells.Item($r, 6).Value = $vals[0]
    $ws.Cells.Item($r, 7).Value = $vals[1]
    $ws.Cells.Item($r, 8).Value = $vals[2]
    $ws.Cells.Item($r, 9).Value = $vals[3]
}

# New row 41
$ws.Range("A41").Formula = "=A40+1"
$ws.Range("B41").Value = 0.025
$ws.Range("C41").Value = 0.025
$ws.Range("D41").Value = 0.0
$ws.Range("E41").Value = 0.0
$ws.Range("F41").Value = 0.0169491525423729
$ws.Range("G41").Value = 0.05
$ws.Range("H41").Value = 0.025
$ws.Range("I41").Value = 0.125

$ws.Range("L6").Select()
